$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet16 ("Sheet16"): selection changes from D7 (sqref A1:J7) to E6 (sqref E6)
# ---------------------------------------------------------------------------
$ws16 = $wb.Worksheets.Item("Sheet16")
$ws16.Activate()
$ws16.Range("E6").Select()

# ---------------------------------------------------------------------------
# Sheet18 ("Sheet18"): selection changes from E27/E27 to (no activeCell)/A1:J7
# ---------------------------------------------------------------------------
$ws18 = $wb.Worksheets.Item("Sheet18")
$ws18.Activate()
$ws18.Range("A1:J7").Select()

# ---------------------------------------------------------------------------
# Sheet23 ("Sheet23"): no longer the active tab (tabSelected removed); the
# selection itself (B4/B4) is unchanged.
# ---------------------------------------------------------------------------
$ws23 = $wb.Worksheets.Item("Sheet23")
$ws23.Range("B4").Select()

# ---------------------------------------------------------------------------
# New sheet "Sheet24" appended at the end of the workbook (copy of Sheet18's
# pre-edit data, with F6 changed from "4--2" to the new string "4--3").
# ---------------------------------------------------------------------------
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws24 = $wb.Worksheets.Add($null, $lastSheet)
$ws24.Name = "Sheet24"

$ws24.Range("A1").Value = "source"
$ws24.Range("B1").Value = "author"
$ws24.Range("C1").Value = "table_name"
$ws24.Range("D1").Value = "dimensions"
$ws24.Range("E1").Value = "dimension_levels_text"
$ws24.Range("F1").Value = "dimension_levels_code"
$ws24.Range("G1").Value = "unit"
$ws24.Range("H1").Value = "interval"
$ws24.Range("I1").Value = "series_name"
$ws24.Range("J1").Value = "table_code"

$ws24.Range("A2").Value = "umar"
$ws24.Range("B2").Value = "mz"
$ws24.Range("C2").Value = "dfg"
$ws24.Range("D2").Value = "dff"
$ws24.Range("E2").Value = 234
$ws24.Range("F2").Value = 1
$ws24.Range("G2").Value = "%"
$ws24.Range("H2").Value = "M"
$ws24.Range("I2").Value = "tesx"
$ws24.Range("J2").Value = "MZ001"

$ws24.Range("A3").Value = "umar"
$ws24.Range("B3").Value = "mz"
$ws24.Range("C3").Value = "dfg"
$ws24.Range("D3").Value = "dff"
$ws24.Range("E3").Value = 1123
$ws24.Range("F3").Value = 3
$ws24.Range("G3").Value = "%"
$ws24.Range("H3").Value = "M"
$ws24.Range("I3").Value = "sdt"
$ws24.Range("J3").Value = "MZ001"

$ws24.Range("A4").Value = "umar"
$ws24.Range("B4").Value = "mz"
$ws24.Range("C4").Value = "dfg"
$ws24.Range("D4").Value = "dff"
$ws24.Range("E4").Value = 1123
$ws24.Range("F4").Value = 2
$ws24.Range("G4").Value = "%"
$ws24.Range("H4").Value = "M"
$ws24.Range("I4").Value = "sdt"
$ws24.Range("J4").Value = "MZ001"

$ws24.Range("A5").Value = "umar"
$ws24.Range("B5").Value = "mz"
$ws24.Range("C5").Value = "hgf"
$ws24.Range("D5").Value = "dim1 -- dim2"
$ws24.Range("E5").Value = "one -- two"
$ws24.Range("F5").Value = "4--2"
$ws24.Range("G5").Value = "%"
$ws24.Range("H5").Value = "A"
$ws24.Range("I5").Value = "sdt"
$ws24.Range("J5").Value = "MZ002"

$ws24.Range("A6").Value = "umar"
$ws24.Range("B6").Value = "mz"
$ws24.Range("C6").Value = "hgf"
$ws24.Range("D6").Value = "dim1 -- dim2"
$ws24.Range("E6").Value = "one -- three"
$ws24.Range("F6").Value = "4--3"
$ws24.Range("G6").Value = "%"
$ws24.Range("H6").Value = "A"
$ws24.Range("I6").Value = "sdt"
$ws24.Range("J6").Value = "MZ002"

$ws24.Range("A7").Value = "umar"
$ws24.Range("B7").Value = "mz"
$ws24.Range("C7").Value = "sdfgs"
$ws24.Range("D7").Value = "dim1 -- dim2-- dim3"
$ws24.Range("E7").Value = "q--e--r"
$ws24.Range("F7").Value = "2--3--4"
$ws24.Range("G7").Value = "%"
$ws24.Range("H7").Value = "A"
$ws24.Range("I7").Value = "sdt"
$ws24.Range("J7").Value = "MZ003"

$ws24.Range("K2").Select()

# ---------------------------------------------------------------------------
# Sheet19 ("Sheet19"): C4/D4 and C5/D5 values change (new shared strings),
# selection changes from I4/A1:K5 to D6/D6, and it becomes the active tab.
# ---------------------------------------------------------------------------
$ws19 = $wb.Worksheets.Item("Sheet19")
$ws19.Activate()
$ws19.Range("C4").Value = "we"
$ws19.Range("D4").Value = "dim1"
$ws19.Range("C5").Value = "we"
$ws19.Range("D5").Value = "dim1"
$ws19.Range("D6").Select()
